# Fruta / hortaliza, semanal
# Insert a new weekly record at row 25, pushing existing rows 25..109 down
# to 26..110 (dimension grows from A1:R109 to A1:R110).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 25 (shifts rows 25..109 down to 26..110)
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new weekly data point
$ws.Range("A25").Value2 = 3
$ws.Range("B25").Value2 = "Femacal de La Calera"
$ws.Range("C25").Value2 = "Coquimbo"
$ws.Range("D25").Value2 = 44560
$ws.Range("E25").Value2 = 5
$ws.Range("F25").Value2 = 100112052
$ws.Range("G25").Value2 = "Albahaca"
$ws.Range("H25").Value2 = "Sin especificar"
$ws.Range("I25").Value2 = "Primera"
$ws.Range("J25").Value2 = 290
$ws.Range("K25").Value2 = 4500
$ws.Range("L25").Value2 = 5000
$ws.Range("M25").Value2 = 4707
$ws.Range("N25").Value2 = "`$/docena de matas"
$ws.Range("O25").Value2 = "Provincia de Quillota"
$ws.Range("P25").Value2 = 784
$ws.Range("Q25").Value2 = 6
$ws.Range("R25").Value2 = "Hortaliza"
